$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Update the "Default System Load" / new "LoadingDetail" columns (row 7 header) ---
$ws.Range("I7").Value = "LoadingDetail"
$ws.Range("H7").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null

# --- Row 8 (PSU830) ---
$ws.Range("G8").Value2 = 0.718
$ws.Range("H8").Value2 = 0.718
$ws.Range("I8").Value = "System (A)"

# --- Row 9 (PSU800) ---
$ws.Range("G9").Value2 = 0.718
$ws.Range("H9").Value2 = 0.718
$ws.Range("I9").Value = "System (A)"

# --- Row 10 (PSU820) ---
$ws.Range("G10").Value2 = 0.568
$ws.Range("H10").Value2 = 0.568
$ws.Range("I10").Value = "System (A)"

# --- Row 11 (PSU821) ---
$ws.Range("G11").Value2 = 0.568
$ws.Range("H11").Value2 = 0.568
$ws.Range("I11").Value = "System (A)"

# Give the new "LoadingDetail" column (I8:I11) the same look as the other
# data cells in that block (style used by B8, D8, E8, ...).
$ws.Range("B8").Copy() | Out-Null
$ws.Range("I8:I11").PasteSpecial(-4122) | Out-Null

# --- Selection, as last edit action performed by the author ---
$ws.Range("G8:G11").Select() | Out-Null
